$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First table cell: "UBND THÀNH PHỐ ĐÀ NẴNG ..." -> "${capHanhChinh} ..."
#    also the paragraph becomes centre-aligned (w:jc w:val="center")
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.Alignment = 1
$p1.Range.Find.Execute('UBND THÀNH PHỐ ĐÀ NẴNG', $true, $false, $false, $false, $false, $true, 1, $false, '${capHanhChinh}', 2) | Out-Null

# ------------------------------------------------------------------
# 2) "THANH TRA THÀNH PHỐ" -> "${coQuanTrucThuoc}"
# ------------------------------------------------------------------
$d.Content.Find.Execute('THANH TRA THÀNH PHỐ', $true, $false, $false, $false, $false, $true, 1, $false, '${coQuanTrucThuoc}', 2) | Out-Null

# ------------------------------------------------------------------
# 3) Append a period after "${nguoiDungDon}" in the first occurrence
#    (" ${nguoiDungDon}" -> " ${nguoiDungDon}.")
# ------------------------------------------------------------------
$d.Content.Find.Execute(' ${nguoiDungDon}', $true, $false, $false, $false, $false, $true, 1, $false, ' ${nguoiDungDon}.', 1) | Out-Null

# ------------------------------------------------------------------
# 4) Append a period after "${diaChi}"
# ------------------------------------------------------------------
$d.Content.Find.Execute(' ${diaChi}', $true, $false, $false, $false, $false, $true, 1, $false, ' ${diaChi}.', 2) | Out-Null

# ------------------------------------------------------------------
# 5) Merge "Nội dung đơn: " + "${noiDung}" into a single run and add a
#    trailing period: "Nội dung đơn: ${noiDung}."
# ------------------------------------------------------------------
$d.Content.Find.Execute('Nội dung đơn: ${noiDung}', $true, $false, $false, $false, $false, $true, 1, $false, 'Nội dung đơn: ${noiDung}.', 2) | Out-Null

# ------------------------------------------------------------------
# 6) Merge ".(" + "3" + ") giải quyết ngày .../..../... (nếu có)" into one
#    run and add a trailing period.
# ------------------------------------------------------------------
$d.Content.Find.Execute('.(3) giải quyết ngày ………/……../……… (nếu có)', $true, $false, $false, $false, $false, $true, 1, $false, '.(3) giải quyết ngày ………/……../……… (nếu có).', 2) | Out-Null

# ------------------------------------------------------------------
# 7) Merge "(" + "3" + ") " into a single run "(3) " in the last
#    paragraph (scoped search avoids touching the other "(3)" above).
# ------------------------------------------------------------------
$p33 = $d.Paragraphs.Item(33)
$p33.Range.Find.Execute('(3) ', $true, $false, $false, $false, $false, $true, 1, $false, '(3) ', 2) | Out-Null

Write-Host "done"
